$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" '65.618.65'
Set-TextValue $ws "E2" '  -0.69%  '
Set-TextValue $ws "D3" '3.442.36'
Set-TextValue $ws "E3" '  -2.61%  '
Set-TextValue $ws "E4" '  -0.13%  '
Set-TextValue $ws "D5" '589.59'
Set-TextValue $ws "E5" '  -1.95%  '
Set-TextValue $ws "D6" '138.23'
Set-TextValue $ws "E6" '  -5.23%  '
Set-TextValue $ws "D7" '3.442.91'
Set-TextValue $ws "E7" '  -2.60%  '
Set-TextValue $ws "E8" '  +0.01%  '
Set-TextValue $ws "D9" '0.499'
Set-TextValue $ws "E9" '  +0.04%  '
Set-TextValue $ws "E10" '  -5.30%  '
Set-TextValue $ws "E11" '  -8.72%  '
Set-TextValue $ws "D12" '0.379'
Set-TextValue $ws "E12" '  -6.86%  '
Set-TextValue $ws "D13" '4.026.42'
Set-TextValue $ws "E13" '  -2.82%  '
Set-TextValue $ws "D14" '0.0000182'
Set-TextValue $ws "E14" '  -9.61%  '
Set-TextValue $ws "D15" '26.51'
Set-TextValue $ws "E15" '  -8.36%  '
Set-TextValue $ws "D16" '3.447.24'
Set-TextValue $ws "E16" '  -2.49%  '
Set-TextValue $ws "D17" '65.576.14'
Set-TextValue $ws "E17" '  -0.83%  '
Set-TextValue $ws "E18" '  -1.55%  '
Set-TextValue $ws "D19" '9.91'
Set-TextValue $ws "E19" '  -9.79%  '
Set-TextValue $ws "E20" '  -5.52%  '
Set-TextValue $ws "D21" '13.74'
Set-TextValue $ws "E21" '  -5.66%  '
Set-TextValue $ws "D22" '393.26'
Set-TextValue $ws "E22" '  -6.28%  '
Set-TextValue $ws "D23" '0.557'
Set-TextValue $ws "E23" '  -7.20%  '
Set-TextValue $ws "D24" '73.41'
Set-TextValue $ws "E24" '  -5.61%  '
Set-TextValue $ws "D25" '1.00'
Set-TextValue $ws "E25" '  -0.06%  '
Set-TextValue $ws "D26" '3.578.53'
Set-TextValue $ws "E26" '  -2.80%  '
Set-TextValue $ws "D27" '0.0000108'
Set-TextValue $ws "E27" '  -7.15%  '
Set-TextValue $ws "E28" '  +0.00%  '
Set-TextValue $ws "D29" '7.23'
Set-TextValue $ws "E29" '  -6.58%  '
Set-TextValue $ws "D30" '8.28'
Set-TextValue $ws "E30" '  -8.81%  '
Set-TextValue $ws "E31" '  -9.00%  '
Set-TextValue $ws "D32" '3.448.11'
Set-TextValue $ws "E32" '  -2.46%  '
Set-TextValue $ws "E33" '  +0.01%  '
Set-TextValue $ws "D34" '0.146'
Set-TextValue $ws "E34" '  -5.67%  '
Set-TextValue $ws "D35" '23.05'
Set-TextValue $ws "E35" '  -5.20%  '
Set-TextValue $ws "D36" '173.03'
Set-TextValue $ws "E36" '  -0.62%  '
Set-TextValue $ws "D37" '6.91'
Set-TextValue $ws "E37" '  -8.56%  '
Set-TextValue $ws "E38" '  -6.83%  '
Set-TextValue $ws "D39" '1.48'
Set-TextValue $ws "E39" '  -6.88%  '
Set-TextValue $ws "E40" '  -7.82%  '
Set-TextValue $ws "D41" '0.0766'
Set-TextValue $ws "E41" '  -6.73%  '
Set-TextValue $ws "D42" '0.826'
Set-TextValue $ws "E42" '  -3.74%  '
Set-TextValue $ws "D43" '43.73'
Set-TextValue $ws "E43" '  -4.04%  '
Set-TextValue $ws "E44" '  -0.08%  '
Set-TextValue $ws "D45" '4.43'
Set-TextValue $ws "E45" '  -12.86%  '
Set-TextValue $ws "D46" '1.63'
Set-TextValue $ws "E46" '  -8.84%  '
Set-TextValue $ws "E47" '  +3.04%  '
Set-TextValue $ws "D48" '22.76'
Set-TextValue $ws "E48" '  +1.03%  '
Set-TextValue $ws "D49" '6.57'
Set-TextValue $ws "E49" '  -7.40%  '
Set-TextValue $ws "E50" '  -12.15%  '
Set-TextValue $ws "D51" '2.214.94'
Set-TextValue $ws "E51" '  -6.73%  '
